$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Process rows top-to-bottom, left-to-right, matching natural save order so
# the shared string table ends up built in the same sequence as the source edit.

# Row 13: mark response column
$ws.Range("F13").Value = "X"

# Row 67: mark response column, update observation
$ws.Range("F67").Value = "X"
$ws.Range("G67").Value = "Respuesta automática"

# Row 82: mark response column, update observation
$ws.Range("F82").Value = "X"
$ws.Range("G82").Value = "No están interesados"

# Row 85: mark response column, update observation
$ws.Range("F85").Value = "X"
$ws.Range("G85").Value = "Escribir correo a: servicliente@panamericana.com.co"

# Row 88: mark response column, update observation, then rename sponsor
$ws.Range("F88").Value = "X"
$ws.Range("G88").Value = "Enviar correo a: talentinquiries@ebat.com"
$ws.Range("C88").Value = "eBay"

# Row 90: mark response column
$ws.Range("F90").Value = "X"

# Row 100: mark response column, update observation
$ws.Range("F100").Value = "X"
$ws.Range("G100").Value = "Enviar la información requerida por Instagram"

# Update sheet selection to reflect latest working position
$ws.Range("G67").Select()
